$d = $word.ActiveDocument
$rng = $d.Content
$rng.Start = 0
$rng.End = 0
$res = $rng.Find.Execute("Rp  ", $true, $false, $false, $false, $false,
                         $true, 0, $false, "", 0)
Write-Output ("Found=" + $res + " Start=" + $rng.Start + " End=" + $rng.End + " Text=[" + $rng.Text + "]")
$rng.Select()
Write-Output ("SelStart=" + $word.Selection.Start + " SelEnd=" + $word.Selection.End)
$word.Selection.Collapse(0)
Write-Output ("After collapse SelStart=" + $word.Selection.Start + " SelEnd=" + $word.Selection.End)
$word.Selection.TypeBackspace()
Write-Output ("After backspace SelStart=" + $word.Selection.Start + " SelEnd=" + $word.Selection.End)
